# Update countries & provincias Spain
#
# This applies the refreshed COVID-19 country data pull: updated case
# counts for several countries, three countries (Albania, Uganda,
# Islas Malvinas) that moved earlier in the table because their figures
# were refreshed, and the "Datos actualizados" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [object[]]$Values
    )
    $n = $Values.Count
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $endCol = [char]([int][char]'A' + $n - 1)
    $ws.Range("A${Row}:${endCol}${Row}").Value = $arr
}

# Header timestamp
Set-Row 1 @("Datos actualizados a 10 de Agosto de 2020 a las 17:42")

# Countries with refreshed totals (name unchanged, row position unchanged)
Set-Row 4   @("Estados Unidos", 5210423, 10979, 2666304, 2378363, 0, 139, 165756)
Set-Row 12  @("Chile", 375044, 1988, 347342, 17563, 0, 62, 10139)
Set-Row 15  @("Reino Unido", 311641, 816, 0, 0, 0, 0, 46574)
Set-Row 19  @("Italia", 250825, 259, 202248, 13368, 0, 4, 35209)
Set-Row 22  @("Alemania", 217329, 48, 197900, 10169, 0, 0, 9260)
Set-Row 38  @("Republica Dominicana", 80499, 767, 44910, 34261, 0, 19, 1328)
Set-Row 47  @("Singapur", 55292, 188, 49609, 5656, 0, 0, 27)
Set-Row 53  @("Barein", 44011, 0, 40967, 2881, 0, 1, 163)
Set-Row 64  @("Moldavia", 27841, 181, 19300, 7691, 0, 5, 850)

# Albania refreshed and moved ahead of Mauritania / Libano (which both
# shift down one row, keeping their own data unchanged)
Set-Row 98  @("Albania", 6536, 125, 3379, 2957, 0, 1, 200)
Set-Row 99  @("Mauritania", 6523, 0, 5527, 839, 0, 0, 157)
Set-Row 100 @("Libano", 6517, 0, 2127, 4312, 0, 0, 78)

# Uganda refreshed and moved ahead of Letonia (which shifts down one row,
# keeping its own data unchanged)
Set-Row 141 @("Uganda", 1297, 14, 1137, 151, 0, 2, 9)
Set-Row 142 @("Letonia", 1290, 0, 1070, 188, 0, 0, 32)

# Santa Lucia / Timor Oriental swap order (figures identical either way)
Set-Row 202 @("Santa Lucia", 25, 0, 24, 1, 0, 0, 0)
Set-Row 203 @("Timor Oriental", 25, 0, 24, 1, 0, 0, 0)

# Islas Malvinas / Montserrat swap order and their figures
Set-Row 213 @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
Set-Row 214 @("Montserrat", 13, 0, 12, 0, 0, 0, 1)
